$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 becomes a real number (was stored as inline string "3524")
$ws.Range("A2").Value = 3524

# New numeric rows
$ws.Range("A3").Value = 2312
$ws.Range("A4").Value = 65543

# New row stored as text (inline string "321"), not as a number
$ws.Range("A5").Value = "'321"
